# Carlos Sahagun resume update:
#   - new contact email
#   - URL text split across a run boundary (around the pre-existing
#     _GoBack bookmark)
#   - three bullets get proofing marks (w:proofErr) around words Word's
#     spell/grammar checker flags, splitting the host run into pieces
#
# Because the `w:proofErr` elements and the mid-run bookmark split have no
# first-class Word object-model property, each affected paragraph is
# rebuilt from its own WordOpenXML (Flat OPC) and written back with
# Range.InsertXML, which replaces exactly the contents of the Range it is
# called on.

$d = $word.ActiveDocument

function Get-ParagraphContaining([string]$needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            return $p
        }
    }
    throw ("Paragraph containing '" + $needle + "' not found")
}

function New-FlatOpcParagraph([string]$paragraphXml) {
    # Wraps a single <w:p>...</w:p> fragment in the minimal Flat OPC
    # envelope Range.InsertXML expects (the same shape Range.WordOpenXML
    # round-trips).
    return (
        '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" ' +
                'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                '<pkg:xmlData>' +
                    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                        '<w:body>' + $paragraphXml + '</w:body>' +
                    '</w:document>' +
                '</pkg:xmlData>' +
            '</pkg:part>' +
        '</pkg:package>'
    )
}

function Set-ParagraphXml([string]$needle, [string]$paragraphXml) {
    $para = Get-ParagraphContaining $needle
    $xml = New-FlatOpcParagraph $paragraphXml
    $para.Range.InsertXML($xml) | Out-Null
    Write-Output ("updated paragraph containing: " + $needle)
}

# --- 1. contact e-mail -------------------------------------------------
$emailPara = @'
<w:p w:rsidR="005B3A33" w:rsidRPr="006268F8" w:rsidRDefault="005B3A33" w:rsidP="00F40A10"><w:pPr><w:jc w:val="right"/><w:rPr><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="00697C6E"><w:rPr><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>sahagun.gino@gmail.com</w:t></w:r></w:p>
'@
Set-ParagraphXml 'sahagunc@spu.edu' $emailPara

# --- 2. personal site URL, split around the _GoBack bookmark -----------
$urlPara = @'
<w:p w:rsidR="005B3A33" w:rsidRPr="006268F8" w:rsidRDefault="00FC3BB4" w:rsidP="00F40A10"><w:pPr><w:jc w:val="right"/><w:rPr><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>https://ginosahagun.githu</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>b.io</w:t></w:r></w:p>
'@
Set-ParagraphXml 'ginosahagun.github.io' $urlPara

# --- 3. "Savour & Sip" bullet, spell-check mark around "Savour" --------
$savourPara = @'
<w:p w:rsidR="005B3A33" w:rsidRPr="00F2756B" w:rsidRDefault="005B3A33" w:rsidP="005B3A33"><w:pPr><w:rPr><w:i/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:i/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Received Commendation for Web App</w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Received 3</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>rd</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> place for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Savour</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> &amp; Sip App in the Social Venture Competition </w:t></w:r><w:r w:rsidRPr="00270697"><w:rPr><w:b/><w:i/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Capable Multi-tasker</w:t></w:r><w:r w:rsidRPr="00270697"><w:rPr><w:i/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>:</w:t></w:r><w:r w:rsidRPr="00F2756B"><w:rPr><w:i/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00F2756B"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Balanced student </w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>work-load</w:t></w:r><w:r w:rsidRPr="00F2756B"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> while working part-time jobs</w:t></w:r><w:r w:rsidRPr="00F2756B"><w:rPr><w:i/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
Set-ParagraphXml 'Savour' $savourPara

# --- 4. SemTech bullet, spell-check mark around "SemTech" --------------
$semtechPara = @'
<w:p w:rsidR="005B3A33" w:rsidRPr="00E80BEC" w:rsidRDefault="005B3A33" w:rsidP="005B3A33"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Created a test-ecosystem (scale-able project) of a cloud, server, and end-device by using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>S</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>emT</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>ech</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> radio devices as a gateway and end-device</w:t></w:r></w:p>
'@
Set-ParagraphXml 'cloud, server' $semtechPara

# --- 5. excel-file bullet, grammar mark around "a" ----------------------
$excelPara = @'
<w:p w:rsidR="005B3A33" w:rsidRPr="00E80BEC" w:rsidRDefault="005B3A33" w:rsidP="005B3A33"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00E80BEC"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Made a desktop through the electron framework</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00E80BEC"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">that converts two files and makes the following changes </w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">in </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> excel file</w:t></w:r></w:p>
'@
Set-ParagraphXml 'in a excel file' $excelPara

Write-Output "done"
